# Daily attendance processing - 2025-12-05 07:02:59
#
# The "Recorded By" column (G) lists the accounts that touched a session's
# attendance record, as a comma-separated string. The recording/export logic
# now lists the "System" entries first (exact "System" before lowercase
# "system"), followed by the remaining account names in their original
# relative order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

$colRecordedBy = 7   # column G

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, $colRecordedBy)
    $text = $cell.Value2

    if ($text -eq $null -or $text -eq "") {
        continue
    }

    $parts = $text -split ", "

    $systemExact = @()
    $systemLower = @()
    $rest = @()

    foreach ($part in $parts) {
        if ($part.Equals("System")) {
            $systemExact += $part
        } elseif ($part.Equals("system")) {
            $systemLower += $part
        } else {
            $rest += $part
        }
    }

    $reordered = $systemExact + $systemLower + $rest
    $newText = $reordered -join ", "

    if ($newText -ne $text) {
        $cell.Value = $newText
    }
}

Write-Output "Recorded By reorder complete: processed rows 2..$lastRow"
